$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply weekly refresh: re-map each data row (2-22) to the updated
# source record (columns D, L:T) per the new fetch ordering.

# Row 2 <- source record (prev row 18)
$ws.Range("D2").Value = 44875
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 400
$ws.Range("N2").Value = 7500
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7750
$ws.Range("Q2").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R2").Value = 'Provincia de Curicó'
$ws.Range("S2").Value = 5167
$ws.Range("T2").Value = 1.5

# Row 3 <- source record (prev row 22)
$ws.Range("D3").Value = 44169
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 400
$ws.Range("N3").Value = 5500
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 5750
$ws.Range("Q3").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R3").Value = 'Provincia de Curicó'
$ws.Range("S3").Value = 3833
$ws.Range("T3").Value = 1.5

# Row 4 <- source record (prev row 21)
$ws.Range("D4").Value = 44166
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6500
$ws.Range("P4").Value = 6250
$ws.Range("Q4").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R4").Value = 'Provincia de Curicó'
$ws.Range("S4").Value = 4167
$ws.Range("T4").Value = 1.5

# Row 5 <- source record (prev row 17)
$ws.Range("D5").Value = 44176
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 5500
$ws.Range("Q5").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R5").Value = 'Provincia de Curicó'
$ws.Range("S5").Value = 3667
$ws.Range("T5").Value = 1.5

# Row 6 <- source record (prev row 15)
$ws.Range("D6").Value = 44516
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 3700
$ws.Range("O6").Value = 3800
$ws.Range("P6").Value = 3750
$ws.Range("Q6").Value = '$/kilo'
$ws.Range("R6").Value = 'Región del Maule'
$ws.Range("S6").Value = 3750
$ws.Range("T6").Value = 1

# Row 7 <- source record (prev row 20)
$ws.Range("D7").Value = 44876
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 7500
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 7750
$ws.Range("Q7").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R7").Value = 'Provincia de Curicó'
$ws.Range("S7").Value = 5167
$ws.Range("T7").Value = 1.5

# Row 8 <- source record (prev row 12)
$ws.Range("D8").Value = 44553
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 400
$ws.Range("N8").Value = 5000
$ws.Range("O8").Value = 5500
$ws.Range("P8").Value = 5250
$ws.Range("Q8").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R8").Value = 'Región del Maule'
$ws.Range("S8").Value = 3500
$ws.Range("T8").Value = 1.5

# Row 9 <- source record (prev row 2)
$ws.Range("D9").Value = 44882
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 7500
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 7750
$ws.Range("Q9").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R9").Value = 'Provincia de Curicó'
$ws.Range("S9").Value = 5167
$ws.Range("T9").Value = 1.5

# Row 10 <- source record (prev row 5)
$ws.Range("D10").Value = 44547
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 400
$ws.Range("N10").Value = 5000
$ws.Range("O10").Value = 5500
$ws.Range("P10").Value = 5250
$ws.Range("Q10").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R10").Value = 'Región del Maule'
$ws.Range("S10").Value = 3500
$ws.Range("T10").Value = 1.5

# Row 11 <- source record (prev row 4)
$ws.Range("D11").Value = 44523
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 3700
$ws.Range("O11").Value = 3800
$ws.Range("P11").Value = 3750
$ws.Range("Q11").Value = '$/kilo'
$ws.Range("R11").Value = 'Región del Maule'
$ws.Range("S11").Value = 3750
$ws.Range("T11").Value = 1

# Row 12 <- source record (prev row 6)
$ws.Range("D12").Value = 44530
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 160
$ws.Range("N12").Value = 3600
$ws.Range("O12").Value = 3700
$ws.Range("P12").Value = 3650
$ws.Range("Q12").Value = '$/kilo'
$ws.Range("R12").Value = 'Región del Maule'
$ws.Range("S12").Value = 3650
$ws.Range("T12").Value = 1

# Row 13 <- source record (prev row 11)
$ws.Range("D13").Value = 44551
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 400
$ws.Range("N13").Value = 5000
$ws.Range("O13").Value = 5500
$ws.Range("P13").Value = 5250
$ws.Range("Q13").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R13").Value = 'Región del Maule'
$ws.Range("S13").Value = 3500
$ws.Range("T13").Value = 1.5

# Row 14 <- source record (prev row 9)
$ws.Range("D14").Value = 44159
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 6500
$ws.Range("O14").Value = 7000
$ws.Range("P14").Value = 6750
$ws.Range("Q14").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R14").Value = 'Provincia de Curicó'
$ws.Range("S14").Value = 4500
$ws.Range("T14").Value = 1.5

# Row 15 <- source record (prev row 13)
$ws.Range("D15").Value = 44519
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 3700
$ws.Range("O15").Value = 3800
$ws.Range("P15").Value = 3750
$ws.Range("Q15").Value = '$/kilo'
$ws.Range("R15").Value = 'Región del Maule'
$ws.Range("S15").Value = 3750
$ws.Range("T15").Value = 1

# Row 16 <- source record (prev row 10)
$ws.Range("D16").Value = 44880
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 7500
$ws.Range("O16").Value = 8000
$ws.Range("P16").Value = 7750
$ws.Range("Q16").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R16").Value = 'Provincia de Curicó'
$ws.Range("S16").Value = 5167
$ws.Range("T16").Value = 1.5

# Row 17 <- source record (prev row 16)
$ws.Range("D17").Value = 44873
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 400
$ws.Range("N17").Value = 7500
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 7750
$ws.Range("Q17").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R17").Value = 'Provincia de Curicó'
$ws.Range("S17").Value = 5167
$ws.Range("T17").Value = 1.5

# Row 18 <- source record (prev row 19)
$ws.Range("D18").Value = 44537
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 400
$ws.Range("N18").Value = 5000
$ws.Range("O18").Value = 5500
$ws.Range("P18").Value = 5250
$ws.Range("Q18").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R18").Value = 'Región del Maule'
$ws.Range("S18").Value = 3500
$ws.Range("T18").Value = 1.5

# Row 19 <- source record (prev row 3)
$ws.Range("D19").Value = 44544
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 400
$ws.Range("N19").Value = 5000
$ws.Range("O19").Value = 5500
$ws.Range("P19").Value = 5250
$ws.Range("Q19").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R19").Value = 'Región del Maule'
$ws.Range("S19").Value = 3500
$ws.Range("T19").Value = 1.5

# Row 20 <- source record (prev row 7)
$ws.Range("D20").Value = 44162
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 7000
$ws.Range("O20").Value = 7000
$ws.Range("P20").Value = 7000
$ws.Range("Q20").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R20").Value = 'Provincia de Curicó'
$ws.Range("S20").Value = 4667
$ws.Range("T20").Value = 1.5

# Row 21 <- source record (prev row 8)
$ws.Range("D21").Value = 44162
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 6500
$ws.Range("O21").Value = 6500
$ws.Range("P21").Value = 6500
$ws.Range("Q21").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R21").Value = 'Provincia de Curicó'
$ws.Range("S21").Value = 4333
$ws.Range("T21").Value = 1.5

# Row 22 <- source record (prev row 14)
$ws.Range("D22").Value = 44533
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 400
$ws.Range("N22").Value = 3500
$ws.Range("O22").Value = 3600
$ws.Range("P22").Value = 3550
$ws.Range("Q22").Value = '$/kilo'
$ws.Range("R22").Value = 'Región del Maule'
$ws.Range("S22").Value = 3550
$ws.Range("T22").Value = 1
